# Fruta / hortaliza, semanal
# The weekly refresh reshuffles which market-day data (Fecha, Volumen,
# Precio minimo/maximo/promedio, Origen, Precio $/Kg) lands on which row.
# Build the row -> source-row mapping (data for row $r becomes the data
# that used to live on row $map[$r]), snapshot every source row first,
# then write the snapshotted values back out so that overlapping reads
# and writes never clobber each other.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    2  = 9
    3  = 22
    4  = 5
    5  = 21
    6  = 6
    7  = 13
    8  = 20
    9  = 11
    10 = 18
    11 = 16
    12 = 2
    13 = 14
    14 = 15
    15 = 25
    16 = 24
    17 = 17
    18 = 7
    19 = 3
    20 = 8
    21 = 12
    22 = 4
    23 = 23
    24 = 10
    25 = 19
}

# Snapshot the columns that move (D, J, K, L, M, O, P) for every data row
# before any writes happen.
$snapshot = @{}
foreach ($r in 2..25) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
    }
}

foreach ($r in 2..25) {
    $src = $snapshot[$map[$r]]
    $ws.Cells.Item($r, 4).Value2 = $src.D
    $ws.Cells.Item($r, 10).Value2 = $src.J
    $ws.Cells.Item($r, 11).Value2 = $src.K
    $ws.Cells.Item($r, 12).Value2 = $src.L
    $ws.Cells.Item($r, 13).Value2 = $src.M
    $ws.Cells.Item($r, 15).Value2 = $src.O
    $ws.Cells.Item($r, 16).Value2 = $src.P
}
